$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Addr, $Text)
    $c = $Sheet.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Text
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "51.850.14"
Set-TextCell $ws "E2" "  +1.29%  "

Set-TextCell $ws "D3" "2.811.56"
Set-TextCell $ws "E3" "  +1.63%  "

Set-TextCell $ws "E4" "  -0.09%  "

Set-TextCell $ws "D5" "351.05"
Set-TextCell $ws "E5" "  -0.81%  "

Set-TextCell $ws "D6" "113.27"
Set-TextCell $ws "E6" "  +4.80%  "

Set-TextCell $ws "D7" "0.559"
Set-TextCell $ws "E7" "  +1.87%  "

Set-TextCell $ws "E8" "  +0.03%  "

Set-TextCell $ws "E9" "  +5.70%  "

Set-TextCell $ws "D10" "40.32"
Set-TextCell $ws "E10" "  +2.10%  "

Set-TextCell $ws "E11" "  -0.88%  "

Set-TextCell $ws "D12" "0.0848"
Set-TextCell $ws "E12" "  +1.76%  "

Set-TextCell $ws "D13" "19.93"
Set-TextCell $ws "E13" "  -0.04%  "

Set-TextCell $ws "D14" "7.81"
Set-TextCell $ws "E14" "  +3.82%  "

Set-TextCell $ws "D15" "3.254.15"
Set-TextCell $ws "E15" "  +1.65%  "

Set-TextCell $ws "D16" "0.969"
Set-TextCell $ws "E16" "  +4.28%  "

Set-TextCell $ws "D17" "2.811.00"
Set-TextCell $ws "E17" "  +1.87%  "

Set-TextCell $ws "D18" "51.855.82"
Set-TextCell $ws "E18" "  +1.36%  "

Set-TextCell $ws "D19" "3.37"
Set-TextCell $ws "E19" "  +9.37%  "

Set-TextCell $ws "E20" "  -0.75%  "

Set-TextCell $ws "D21" "13.54"
Set-TextCell $ws "E21" "  +3.18%  "

Set-TextCell $ws "D22" "0.0₃0977"
Set-TextCell $ws "E22" "  +1.67%  "

Set-TextCell $ws "D23" "70.63"
Set-TextCell $ws "E23" "  +1.42%  "

Set-TextCell $ws "D24" "268.59"
Set-TextCell $ws "E24" "  +1.38%  "

Set-TextCell $ws "E25" "  +1.79%  "

Set-TextCell $ws "D26" "26.25"
Set-TextCell $ws "E26" "  +1.15%  "

Set-TextCell $ws "D27" "1.00"
Set-TextCell $ws "E27" "  -0.11%  "

Set-TextCell $ws "D28" "0.163"
Set-TextCell $ws "E28" "  +0.63%  "

Set-TextCell $ws "B29" "InjectiveProtocol"
Set-TextCell $ws "C29" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D29" "38.71"
Set-TextCell $ws "E29" "  +10.98%  "

Set-TextCell $ws "B30" "Cosmos"
Set-TextCell $ws "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D30" "10.49"
Set-TextCell $ws "E30" "  +3.13%  "

Set-TextCell $ws "E31" "  +0.91%  "

Set-TextCell $ws "D32" "52.78"
Set-TextCell $ws "E32" "  +1.84%  "

Set-TextCell $ws "D33" "6.17"
Set-TextCell $ws "E33" "  +1.57%  "

Set-TextCell $ws "B34" "RenderToken"
Set-TextCell $ws "C34" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D34" "5.70"
Set-TextCell $ws "E34" "  +3.42%  "

Set-TextCell $ws "B35" "Hedera"
Set-TextCell $ws "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D35" "0.0902"
Set-TextCell $ws "E35" "  +8.53%  "

Set-TextCell $ws "E36" "  +2.75%  "

Set-TextCell $ws "D37" "0.999"
Set-TextCell $ws "E37" "  -0.20%  "

Set-TextCell $ws "D38" "19.01"
Set-TextCell $ws "E38" "  +4.70%  "

Set-TextCell $ws "E39" "  +2.30%  "

Set-TextCell $ws "D40" "2.01"
Set-TextCell $ws "E40" "  +3.05%  "

Set-TextCell $ws "D41" "2.57"
Set-TextCell $ws "E41" "  +2.51%  "

Set-TextCell $ws "E42" "  +1.99%  "

Set-TextCell $ws "B43" "WEMIXToken"
Set-TextCell $ws "C43" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D43" "2.24"
Set-TextCell $ws "E43" "  +2.06%  "

Set-TextCell $ws "B44" "EnergySwap"
Set-TextCell $ws "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D44" "22.24"
Set-TextCell $ws "E44" "  +1.19%  "

Set-TextCell $ws "D45" "120.85"
Set-TextCell $ws "E45" "  +0.25%  "

Set-TextCell $ws "E46" "  +8.30%  "

Set-TextCell $ws "E47" "  +10.30%  "

Set-TextCell $ws "D48" "2.145.58"
Set-TextCell $ws "E48" "  +2.80%  "

Set-TextCell $ws "D49" "1.01"
Set-TextCell $ws "E49" "  +10.33%  "

Set-TextCell $ws "D50" "0.227"
Set-TextCell $ws "E50" "  +19.86%  "

Set-TextCell $ws "B51" "THORChain"
Set-TextCell $ws "C51" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws "D51" "5.52"
Set-TextCell $ws "E51" "  +0.36%  "

